$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new data rows (image-transfer test results) ---
# Shared-string / row ordering follows the original author's edit sequence:
# first the "Local_to_PI_Redhawk_Logo" rows (58-63), then the
# "Remote_from_PI_Peter_Jamieson" rows (40-45) filling the earlier gap,
# then the final "Remote_from_PI_Redhawk_Logo" rows (64-69).

$ws.Range("A58").Value = "PC_Local_to_PI_Redhawk_Logo_double_bmp"
$ws.Range("A58").Font.Bold = $true
$ws.Range("B58").Value = 3.6657
$ws.Range("C58").Value = 3.5435
$ws.Range("D58").Formula = "=B58-C58"

$ws.Range("A59").Value = "PC_Local_to_PI_Redhawk_Logo_double_png"
$ws.Range("A59").Font.Bold = $true
$ws.Range("B59").Value = 3.5545
$ws.Range("C59").Value = 3.3877
$ws.Range("D59").Formula = "=B59-C59"

$ws.Range("A60").Value = "PC_Local_to_PI_Redhawk_Logo_full_bmp"
$ws.Range("A60").Font.Bold = $true
$ws.Range("B60").Value = 3.5705
$ws.Range("C60").Value = 3.4543
$ws.Range("D60").Formula = "=B60-C60"

$ws.Range("A61").Value = "PC_Local_to_PI_Redhawk_Logo_full_png"
$ws.Range("A61").Font.Bold = $true
$ws.Range("B61").Value = 3.5395
$ws.Range("C61").Value = 3.3852
$ws.Range("D61").Formula = "=B61-C61"

$ws.Range("A62").Value = "PC_Local_to_PI_Redhawk_Logo_half_bmp"
$ws.Range("A62").Font.Bold = $true
$ws.Range("B62").Value = 3.5925
$ws.Range("C62").Value = 3.4387
$ws.Range("D62").Formula = "=B62-C62"

$ws.Range("A63").Value = "PC_Local_to_PI_Redhawk_Logo_half_png"
$ws.Range("A63").Font.Bold = $true
$ws.Range("B63").Value = 3.5673
$ws.Range("C63").Value = 3.3717
$ws.Range("D63").Formula = "=B63-C63"

$ws.Range("A40").Value = "PC_Remote_from_PI_Peter_Jamieson_Double_bmp"
$ws.Range("A40").Font.Bold = $true
$ws.Range("B40").Value = 3.6937
$ws.Range("C40").Value = 3.6193
$ws.Range("D40").Formula = "=B40-C40"

$ws.Range("A41").Value = "PC_Remote_from_PI_Peter_Jamieson_Double_png"
$ws.Range("A41").Font.Bold = $true
$ws.Range("B41").Value = 3.6281
$ws.Range("C41").Value = 3.53
$ws.Range("D41").Formula = "=B41-C41"

$ws.Range("A42").Value = "PC_Remote_from_PI_Peter_Jamieson_full_bmp"
$ws.Range("A42").Font.Bold = $true
$ws.Range("B42").Value = 3.6173
$ws.Range("C42").Value = 3.5208
$ws.Range("D42").Formula = "=B42-C42"

$ws.Range("A43").Value = "PC_Remote_from_PI_Peter_Jamieson_full_png"
$ws.Range("A43").Font.Bold = $true
$ws.Range("B43").Value = 3.5758
$ws.Range("C43").Value = 3.4148
$ws.Range("D43").Formula = "=B43-C43"

$ws.Range("A44").Value = "PC_Remote_from_PI_Peter_Jamieson_half_bmp"
$ws.Range("A44").Font.Bold = $true
$ws.Range("B44").Value = 3.5489
$ws.Range("C44").Value = 3.3899
$ws.Range("D44").Formula = "=B44-C44"

$ws.Range("A45").Value = "PC_Remote_from_PI_Peter_Jamieson_half_png"
$ws.Range("A45").Font.Bold = $true
$ws.Range("B45").Value = 3.5732
$ws.Range("C45").Value = 3.3988
$ws.Range("D45").Formula = "=B45-C45"

$ws.Range("A64").Value = "PC_Remote_from_PI_Redhawk_Logo_double_bmp"
$ws.Range("A64").Font.Bold = $true
$ws.Range("B64").Value = 3.6579
$ws.Range("C64").Value = 3.5529
$ws.Range("D64").Formula = "=B64-C64"

$ws.Range("A65").Value = "PC_Remote_from_PI_Redhawk_Logo_double_png"
$ws.Range("A65").Font.Bold = $true
$ws.Range("B65").Value = 3.5653
$ws.Range("C65").Value = 3.3624
$ws.Range("D65").Formula = "=B65-C65"

$ws.Range("A66").Value = "PC_Remote_from_PI_Redhawk_Logo_full_bmp"
$ws.Range("A66").Font.Bold = $true
$ws.Range("B66").Value = 3.5749
$ws.Range("C66").Value = 3.3888
$ws.Range("D66").Formula = "=B66-C66"

$ws.Range("A67").Value = "PC_Remote_from_PI_Redhawk_Logo_full_png"
$ws.Range("A67").Font.Bold = $true
$ws.Range("B67").Value = 3.5333
$ws.Range("C67").Value = 3.3941
$ws.Range("D67").Formula = "=B67-C67"

$ws.Range("A68").Value = "PC_Remote_from_PI_Redhawk_Logo_half_bmp"
$ws.Range("A68").Font.Bold = $true
$ws.Range("B68").Value = 3.5323
$ws.Range("C68").Value = 3.3715
$ws.Range("D68").Formula = "=B68-C68"

$ws.Range("A69").Value = "PC_Remote_from_PI_Redhawk_Logo_half_png"
$ws.Range("A69").Font.Bold = $true
$ws.Range("B69").Value = 3.527
$ws.Range("C69").Value = 3.3954
$ws.Range("D69").Formula = "=B69-C69"

# --- Update selection / view state to match the saved workbook ---
$ws.Range("D69").Select()
